# Module 8 Exercise 1 edits
# ---------------------------------------------------------------
# 1. Remove the "npm install @angular-devkit/core" paragraph.
# 2. In " In "app-routing.module.ts" delete all code ..." change
#    "app-" to "about-" (split into separate runs).
# 3. Fix "Bow we have a basic setup..." -> "Now we have a basic setup..."
#    and move the "_GoBack" bookmark from the end of the document to
#    right after the initial "N".
# ---------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Delete the "npm install @angular-devkit/core" paragraph -------
$rng = $d.Content
$rng.Find.Execute("npm install @angular-devkit/core", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    # Extend the range to swallow the paragraph mark that ends this
    # paragraph so the whole paragraph (incl. its own pilcrow) disappears
    # and the following paragraph's text slides up into this slot.
    $delRange = $d.Range($rng.Start, $rng.End + 1)
    $delRange.Delete()
}

# --- 2. " In "app-routing.module.ts" -> " In "about-routing.module.ts" -
$rng2 = $d.Content
$quote = [char]0x201C
$searchText = " In " + $quote + "app-"
$rng2.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng2.Find.Found) {
    $matchEnd = $rng2.End
    $appStart = $matchEnd - 4

    # Replace "app-" with "about-"
    $appRange = $d.Range($appStart, $matchEnd)
    $appRange.Text = "about-"

    # Force "about" and "-" to materialize as their own separate runs
    # (matching the target markup) by touching formatting on each span.
    $aboutRange = $d.Range($appStart, $appStart + 5)
    $aboutRange.Font.Bold = 1
    $aboutRange.Font.Bold = 0

    $dashRange = $d.Range($appStart + 5, $appStart + 6)
    $dashRange.Font.Bold = 1
    $dashRange.Font.Bold = 0
}

# --- 3. "Bow we have a basic setup..." -> "Now we have a basic setup..."
#        and relocate the "_GoBack" bookmark ----------------------------

# Remove the old "_GoBack" bookmark sitting at the end of the document.
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

$rng3 = $d.Content
$rng3.Find.Execute("Bow we have a basic setup", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng3.Find.Found) {
    $bStart = $rng3.Start

    # Fix the typo: "Bow" -> "Now"
    $bRange = $d.Range($bStart, $bStart + 1)
    $bRange.Text = "N"

    # Re-add the "_GoBack" bookmark right after the "N"
    $afterN = $bStart + 1
    $bmRange = $d.Range($afterN, $afterN)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}
